$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the "through" date in its title (June 15 -> June 16)
$ws.Name = "Through 2022-06-16"

# Update the header text for column B (June 2022 month-to-date label)
$ws.Range("B1").Value = "June 2022 (through June 16)"

# Add new carjacking counts for 2022-06-24 (one additional incident recorded
# per neighborhood/month column)
$ws.Range("H2").Value = 6    # Englewood, June 2021
$ws.Range("AL2").Value = 2   # Englewood, June 2016

$ws.Range("B3").Value = 6    # Auburn Gresham, June 2022
$ws.Range("T3").Value = 3    # Auburn Gresham, June 2019

$ws.Range("T5").Value = 1    # South Shore, June 2019
$ws.Range("AR5").Value = 1   # South Shore, June 2015

$ws.Range("B7").Value = 1    # Grand Boulevard, June 2022

$ws.Range("Z9").Value = 2    # Grand Crossing, June 2018
$ws.Range("AF9").Value = 2   # Grand Crossing, June 2017

$ws.Range("B12").Value = 3   # Roseland, June 2022

$ws.Range("B14").Value = 4   # Austin, June 2022
$ws.Range("AF14").Value = 2  # Austin, June 2017

$ws.Range("N20").Value = 1   # Hyde Park, June 2020

$ws.Range("AL28").Value = 1  # Ukrainian Village, June 2016

$ws.Range("N46").Value = 1   # Avondale, June 2020

$ws.Range("N47").Value = 1   # Belmont Cragin, June 2020

$ws.Range("T85").Value = 1   # River North, June 2019
